$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2569729462136934
$ws.Range("C2").Value = 0.05177388402117344
$ws.Range("D2").Value = 0.07886614501170186
$ws.Range("E2").Value = 0.1679569557938976
$ws.Range("G2").Value = 0.486907209969111
$ws.Range("H2").Value = 0.6345743011770892
$ws.Range("I2").Value = 0.5079287286296754
$ws.Range("K2").Value = 0.2764099386066334
$ws.Range("M2").Value = 0.2175147502625947
$ws.Range("O2").Value = 2.204415478390189
# Row 3
$ws.Range("B3").Value = 0.2247108474927586
$ws.Range("C3").Value = 0.04674439375598638
$ws.Range("D3").Value = 0.07149884586686994
$ws.Range("E3").Value = 0.1572382192002379
$ws.Range("G3").Value = 0.4878594031484127
$ws.Range("H3").Value = 0.6387970036271469
$ws.Range("I3").Value = 0.5127568242317651
$ws.Range("K3").Value = 0.2411291852741329
$ws.Range("M3").Value = 0.1947965997721681
$ws.Range("O3").Value = 2.215010613189818
# Row 4
$ws.Range("B4").Value = 0.2048533372612553
$ws.Range("C4").Value = 0.04363508242363423
$ws.Range("D4").Value = 0.0670079498903533
$ws.Range("E4").Value = 0.1507792560860608
$ws.Range("G4").Value = 0.4888031866333122
$ws.Range("H4").Value = 0.6416829544633629
$ws.Range("I4").Value = 0.5160172512960166
$ws.Range("K4").Value = 0.2193877397840822
$ws.Range("M4").Value = 0.1809048839706691
$ws.Range("O4").Value = 2.222881997568834
# Row 5
$ws.Range("B5").Value = 0.1967495430767769
$ws.Range("C5").Value = 0.04236273944333391
$ws.Range("D5").Value = 0.06518611505656224
$ws.Range("E5").Value = 0.1481777308007111
$ws.Range("G5").Value = 0.4892779692753919
$ws.Range("H5").Value = 0.6429327533882514
$ws.Range("I5").Value = 0.5174202739765157
$ws.Range("K5").Value = 0.2105085712020696
$ws.Range("M5").Value = 0.1752583178814504
$ws.Range("O5").Value = 2.226432906976896
# Row 6
$ws.Range("B6").Value = 0.1954032220679665
$ws.Range("C6").Value = 0.04215115094717703
$ws.Range("D6").Value = 0.06488409944158491
$ws.Range("E6").Value = 0.1477475885074924
$ws.Range("G6").Value = 0.4893622497879448
$ws.Range("H6").Value = 0.6431447366214158
$ws.Range("I6").Value = 0.5176577352977887
$ws.Range("K6").Value = 0.2090330376810527
$ws.Range("M6").Value = 0.1743215811240191
$ws.Range("O6").Value = 2.227043258256344
# Row 7
$ws.Range("B7").Value = 0.2047440933758367
$ws.Range("C7").Value = 0.04361794444285749
$ws.Range("D7").Value = 0.06698334657801297
$ws.Range("E7").Value = 0.1507440475822861
$ws.Range("G7").Value = 0.4888092247347089
$ws.Range("H7").Value = 0.6416995110444077
$ws.Range("I7").Value = 0.5160358718855953
$ws.Range("K7").Value = 0.2192680699170637
$ws.Range("M7").Value = 0.1808286740051699
$ws.Range("O7").Value = 2.222928496835252
# Row 8
$ws.Range("B8").Value = 0.245859311855213
$ws.Range("C8").Value = 0.05004414309478022
$ws.Range("D8").Value = 0.07631913854253014
$ws.Range("E8").Value = 0.1642355987141642
$ws.Range("G8").Value = 0.4871609125308325
$ws.Range("H8").Value = 0.6359694507786742
$ws.Range("I8").Value = 0.5095320009497897
$ws.Range("K8").Value = 0.2642617989471887
$ws.Range("M8").Value = 0.2096696322514831
$ws.Range("O8").Value = 2.207785048234385
# Row 9
$ws.Range("B9").Value = 0.3260842191237998
$ws.Range("C9").Value = 0.06247603981634597
$ws.Range("D9").Value = 0.09488533776642782
$ws.Range("E9").Value = 0.1916749387630574
$ws.Range("G9").Value = 0.4867843483849867
$ws.Range("H9").Value = 0.6270583157499203
$ws.Range("I9").Value = 0.4991284894375063
$ws.Range("K9").Value = 0.3518507282082339
$ws.Range("M9").Value = 0.2666841452789868
$ws.Range("O9").Value = 2.188938458138495
# Row 10
$ws.Range("B10").Value = 0.3847626671695821
$ws.Range("C10").Value = 0.07150467058573895
$ws.Range("D10").Value = 0.1086845880916627
$ws.Range("E10").Value = 0.2124514107873665
$ws.Range("G10").Value = 0.4882581497207923
$ws.Range("H10").Value = 0.6219281191255988
$ws.Range("I10").Value = 0.4929212651769959
$ws.Range("K10").Value = 0.415792659496077
$ws.Range("M10").Value = 0.3088596107593702
$ws.Range("O10").Value = 2.181724855251844
# Row 11
$ws.Range("B11").Value = 0.4113965493761782
$ws.Range("C11").Value = 0.07558894091144452
$ws.Range("D11").Value = 0.1149969510842368
$ws.Range("E11").Value = 0.2220410569560585
$ws.Range("G11").Value = 0.489310810905593
$ws.Range("H11").Value = 0.6199017553545474
$ws.Range("I11").Value = 0.4904100213860048
$ws.Range("K11").Value = 0.4447892809808138
$ws.Range("M11").Value = 0.328110595070882
$ws.Range("O11").Value = 2.179887697615385
# Row 12
$ws.Range("B12").Value = 0.4214731766791999
$ws.Range("C12").Value = 0.07713221301494855
$ws.Range("D12").Value = 0.1173923016344958
$ws.Range("E12").Value = 0.2256925695778307
$ws.Range("G12").Value = 0.4897645407222626
$ws.Range("H12").Value = 0.6191786123671648
$ws.Range("I12").Value = 0.4895040690034484
$ws.Range("K12").Value = 0.4557560659010846
$ws.Range("M12").Value = 0.3354098721920522
$ws.Range("O12").Value = 2.17940000763096
# Row 13
$ws.Range("B13").Value = 0.4193034053889448
$ws.Range("C13").Value = 0.07679999179964625
$ws.Range("D13").Value = 0.1168761984863522
$ws.Range("E13").Value = 0.2249052524076163
$ws.Range("G13").Value = 0.4896643685330559
$ws.Range("H13").Value = 0.6193323885540707
$ws.Range("I13").Value = 0.4896971794973837
$ws.Range("K13").Value = 0.4533947862261698
$ws.Range("M13").Value = 0.3338374268958333
$ws.Range("O13").Value = 2.179495784599482
# Row 14
$ws.Range("B14").Value = 0.412225743875922
$ws.Range("C14").Value = 0.07571597436424327
$ws.Range("D14").Value = 0.1151939181230546
$ws.Range("E14").Value = 0.2223410645058905
$ws.Range("G14").Value = 0.4893470342088762
$ws.Range("H14").Value = 0.6198413762005686
$ws.Range("I14").Value = 0.4903345858442378
$ws.Range("K14").Value = 0.4456918007381319
$ws.Range("M14").Value = 0.3287109237874475
$ws.Range("O14").Value = 2.179843404653013
# Row 15
$ws.Range("B15").Value = 0.4078892781761567
$ws.Range("C15").Value = 0.07505154403561676
$ws.Range("D15").Value = 0.1141641222277769
$ws.Range("E15").Value = 0.2207730526686049
$ws.Range("G15").Value = 0.4891598390472183
$ws.Range("H15").Value = 0.6201589014362838
$ws.Range("I15").Value = 0.4907308782075503
$ws.Range("K15").Value = 0.4409717118503522
$ws.Range("M15").Value = 0.3255720092333618
$ws.Range("O15").Value = 2.180083429005464
# Row 16
$ws.Range("B16").Value = 0.3830208431705557
$ws.Range("C16").Value = 0.07123728802540086
$ws.Range("D16").Value = 0.108272761560599
$ws.Range("E16").Value = 0.2118275080397467
$ws.Range("G16").Value = 0.4881970596891989
$ws.Range("H16").Value = 0.6220667317654147
$ws.Range("I16").Value = 0.4930916742780056
$ws.Range("K16").Value = 0.4138957872912954
$ws.Range("M16").Value = 0.307602822855344
$ws.Range("O16").Value = 2.181874008005565
# Row 17
$ws.Range("B17").Value = 0.3677493155168747
$ws.Range("C17").Value = 0.06889145801925167
$ws.Range("D17").Value = 0.1046675429461175
$ws.Range("E17").Value = 0.2063752987333416
$ws.Range("G17").Value = 0.4877044266618498
$ws.Range("H17").Value = 0.6233158497533964
$ws.Range("I17").Value = 0.4946200234376938
$ws.Range("K17").Value = 0.3972619240363997
$ws.Range("M17").Value = 0.2965959809590473
$ws.Range("O17").Value = 2.183342620630185
# Row 18
$ws.Range("B18").Value = 0.3589599923558637
$ws.Range("C18").Value = 0.06754004681535264
$ws.Range("D18").Value = 0.1025972147187844
$ws.Range("E18").Value = 0.2032523381525024
$ws.Range("G18").Value = 0.4874570449544819
$ws.Range("H18").Value = 0.624063242522638
$ws.Range("I18").Value = 0.4955284948525005
$ws.Range("K18").Value = 0.3876860434157265
$ws.Range("M18").Value = 0.2902712705045971
$ws.Range("O18").Value = 2.184323258271149
# Row 19
$ws.Range("B19").Value = 0.35598314241048
$ws.Range("C19").Value = 0.06708211444295387
$ws.Range("D19").Value = 0.1018968048038147
$ws.Range("E19").Value = 0.2021971841177432
$ws.Range("G19").Value = 0.4873794587702918
$ws.Range("H19").Value = 0.624321266326362
$ws.Range("I19").Value = 0.4958411359142296
$ws.Range("K19").Value = 0.3844423674470931
$ws.Range("M19").Value = 0.2881308863786103
$ws.Range("O19").Value = 2.18467862216184
# Row 20
$ws.Range("B20").Value = 0.3693755733549153
$ws.Range("C20").Value = 0.06914139876708703
$ws.Range("D20").Value = 0.1050509835586979
$ws.Range("E20").Value = 0.2069543479777423
$ws.Range("G20").Value = 0.4877531446904158
$ws.Range("H20").Value = 0.6231798845760039
$ws.Range("I20").Value = 0.4944542841982518
$ws.Range("K20").Value = 0.3990335131150857
$ws.Range("M20").Value = 0.2977670433940602
$ws.Range("O20").Value = 2.183172213782768
# Row 21
$ws.Range("B21").Value = 0.4143048743235909
$ws.Range("C21").Value = 0.07603446802042413
$ws.Range("D21").Value = 0.1156879093809664
$ws.Range("E21").Value = 0.2230936807126085
$ws.Range("G21").Value = 0.4894387461603031
$ws.Range("H21").Value = 0.6196906747367876
$ws.Range("I21").Value = 0.4901461422689302
$ws.Range("K21").Value = 0.4479547278428413
$ws.Range("M21").Value = 0.3302164485146832
$ws.Range("O21").Value = 2.179735652675419
# Row 22
$ws.Range("B22").Value = 0.4436158002302477
$ws.Range("C22").Value = 0.08051991930463487
$ws.Range("D22").Value = 0.122668855138599
$ws.Range("E22").Value = 0.2337590595468839
$ws.Range("G22").Value = 0.4908616741497411
$ws.Range("H22").Value = 0.6176678860780669
$ws.Range("I22").Value = 0.4875928562143059
$ws.Range("K22").Value = 0.4798480137448564
$ws.Range("M22").Value = 0.3514784827137234
$ws.Range("O22").Value = 2.178702195664101
# Row 23
$ws.Range("B23").Value = 0.4279770279260049
$ws.Range("C23").Value = 0.07812775960762508
$ws.Range("D23").Value = 0.118940343584427
$ws.Range("E23").Value = 0.2280559285502548
$ws.Range("G23").Value = 0.4900727843678681
$ws.Range("H23").Value = 0.6187239166048641
$ws.Range("I23").Value = 0.4889315661786462
$ws.Range("K23").Value = 0.4628334232025395
$ws.Range("M23").Value = 0.3401255579831357
$ws.Range("O23").Value = 2.179142723339083
# Row 24
$ws.Range("B24").Value = 0.3686403717959195
$ws.Range("C24").Value = 0.06902840915462605
$ws.Range("D24").Value = 0.1048776227047767
$ws.Range("E24").Value = 0.2066925237414168
$ws.Range("G24").Value = 0.4877310076420542
$ws.Range("H24").Value = 0.6232412632931101
$ws.Range("I24").Value = 0.4945291220779104
$ws.Range("K24").Value = 0.3982326176828224
$ws.Range("M24").Value = 0.2972375958844822
$ws.Range("O24").Value = 2.183248830131902
# Row 25
$ws.Range("B25").Value = 0.3044261292567398
$ws.Range("C25").Value = 0.05913122668015092
$ws.Range("D25").Value = 0.08983489048006277
$ws.Range("E25").Value = 0.1841448455633312
$ws.Range("G25").Value = 0.4865795031560083
$ws.Range("H25").Value = 0.629220144438591
$ws.Range("I25").Value = 0.501690932809165
$ws.Range("K25").Value = 0.3282263025243424
$ws.Range("M25").Value = 0.2512104038531788
$ws.Range("O25").Value = 2.192873487003851
